$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Add new row 31 data to Sheet1 (new LeetCode entry: "Single Number") ---

# Add the hyperlink for B31 first (so the "Single Number" shared string is
# interned before the submission URL, matching natural entry order), matching
# the style used by the other problem-name cells (reuse the existing
# "Hyperlink" style rather than letting Excel mint a fresh one), and set the
# display text to the problem name.
$ws1.Hyperlinks.Add($ws1.Range("B31"), "https://leetcode.com/problems/single-number/", [Type]::Missing, [Type]::Missing, "https://leetcode.com/problems/single-number/")
$ws1.Range("B31").Value = "Single Number"
$ws1.Range("B31").Style = $ws1.Range("B30").Style

$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = 1
$ws1.Range("E31").Value = 121
$ws1.Range("F31").Value = 0.26
$ws1.Range("G31").Value = 18.85
$ws1.Range("H31").Value = 0.14
$ws1.Range("I31").Value = "https://leetcode.com/problems/single-number/submissions/"

# --- Update the active selections on each sheet ---
$ws1.Activate() | Out-Null
$ws1.Range("B35").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("H8").Select() | Out-Null

$ws1.Activate() | Out-Null
